$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the record count values in column D
$ws.Range("D2").Value = 4
$ws.Range("D3").Value = 4

# Correct the CURP value in AB3 (one character shorter)
$ws.Range("AB3").Value = "CAAM970828HCALKR22"

# Widen column AB (CURP) to fit the values better
$ws.Columns("AB").ColumnWidth = 26.67

# Add the new "Posicion codigo" column (AS) at the end of the header row,
# matching the header formatting used by the rest of the nomina block (AE1:AR1)
$ws.Range("AS1").Value = "Posicion codigo"
$ws.Range("AR1").Copy()
$ws.Range("AS1").PasteSpecial(-4122)

$ws.Range("AS2").Value = 27677
$ws.Range("AS3").Value = 2531932

# Move the active selection to the newly added cell, like the author did
$ws.Range("AS2").Select()
